$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 'Augment 100ml PFS'
$ws.Range("E2").Value = '1''s'
$ws.Range("F2").Value = 233
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 'Augment 375 Tablet 18''s'
$ws.Range("E3").Value = '18''s'
$ws.Range("F3").Value = 53
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 'Augment 1.2g IV Injection 1''s'
$ws.Range("E4").Value = '1''s'
$ws.Range("F4").Value = 274
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 'Augment 1gm Tablet 12''s'
$ws.Range("E5").Value = '12''s'
$ws.Range("F5").Value = 78
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 'Augment 625 Tablet 18''s'
$ws.Range("E6").Value = '18''s'
$ws.Range("F6").Value = 342
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 'Biltin 20mg Tablet 20''s'
$ws.Range("E7").Value = '20''s'
$ws.Range("F7").Value = 2009
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 'Bonflex FC Tablet'
$ws.Range("E8").Value = '40 ''s'
$ws.Range("F8").Value = 23
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 'Desodin 5mg Tablet'
$ws.Range("E9").Value = '50 ''s'
$ws.Range("F9").Value = 105
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 'Dinafex 60mg FC Tablet 40''s'
$ws.Range("E10").Value = '40''s'
$ws.Range("F10").Value = 86
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 'Dinafex 50ml Suspension'
$ws.Range("E11").Value = '50 ml'
$ws.Range("F11").Value = 1792
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 'Dinafex 180mg FC Tablet 40''s'
$ws.Range("E12").Value = '40''s'
$ws.Range("F12").Value = 89
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 'Dinafex 120mg FC Tablet 40''s'
$ws.Range("E13").Value = '40''s'
$ws.Range("F13").Value = 688
$ws.Range("C14").Value = 13
$ws.Range("D14").Value = 'Dorenta 100ml Syrup'
$ws.Range("E14").Value = '100 ml'
$ws.Range("F14").Value = 1306
$ws.Range("C15").Value = 14
$ws.Range("D15").Value = 'Etorix 60mg Tablet 50''s'
$ws.Range("E15").Value = '50''s'
$ws.Range("F15").Value = 1389
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 'Etorix 90mg Tablet 40''s'
$ws.Range("E16").Value = '40''s'
$ws.Range("F16").Value = 2721
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 'Etorix 120mg Tablet 30''s'
$ws.Range("E17").Value = '30''s'
$ws.Range("F17").Value = 2839
$ws.Range("C18").Value = 17
$ws.Range("D18").Value = 'Fenobac 10mg Tablet'
$ws.Range("E18").Value = '30 ''s'
$ws.Range("F18").Value = 91
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 'Fenobac 5mg Tablet'
$ws.Range("E19").Value = '50 ''s'
$ws.Range("F19").Value = 35
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 'Flucloxin 250mg Capsule'
$ws.Range("E20").Value = '100 ''s'
$ws.Range("F20").Value = 383
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 'Flucloxin 500mg Capsule 40''s'
$ws.Range("E21").Value = '40''s'
$ws.Range("F21").Value = 3376
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 'Flucloxin 100ml Dry Suspension'
$ws.Range("E22").Value = '100 ml'
$ws.Range("F22").Value = 6195
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 'Flucloxin 500mg IM/IV Injection'
$ws.Range("E23").Value = '1''s'
$ws.Range("F23").Value = 670
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 'Ketonic 60mg IM Injection'
$ws.Range("E24").Value = '1 ''s'
$ws.Range("F24").Value = 1071
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 'Ketonic 10mg Tablet - 30''s'
$ws.Range("E25").Value = '30''s'
$ws.Range("F25").Value = 219
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 'Ketonic 30mg IM/IV Injection'
$ws.Range("E26").Value = '1 ''s'
$ws.Range("F26").Value = 1938
$ws.Range("C27").Value = 26
$ws.Range("D27").Value = 'Levomax 500mg Tablet - 20''s'
$ws.Range("E27").Value = '20 ''s'
$ws.Range("F27").Value = 64
$ws.Range("C28").Value = 27
$ws.Range("D28").Value = 'Levomax 750mg Tablet - 10''s'
$ws.Range("E28").Value = '10 ''s'
$ws.Range("F28").Value = 3
$ws.Range("C29").Value = 28
$ws.Range("D29").Value = 'Lindamax Plus 10gm Gel'
$ws.Range("E29").Value = '10gm'
$ws.Range("F29").Value = 87
$ws.Range("C30").Value = 29
$ws.Range("D30").Value = 'Lindamax 150mg Capsule'
$ws.Range("E30").Value = '30''s'
$ws.Range("F30").Value = 11
$ws.Range("C31").Value = 30
$ws.Range("D31").Value = 'Lindamax 25ml Lotion'
$ws.Range("E31").Value = '1''s'
$ws.Range("F31").Value = 59
$ws.Range("C32").Value = 31
$ws.Range("D32").Value = 'Lindamax 300mg Capsule'
$ws.Range("E32").Value = '30''s'
$ws.Range("F32").Value = 54
$ws.Range("C33").Value = 32
$ws.Range("D33").Value = 'Mebidal Tablet'
$ws.Range("E33").Value = '200 ''s'
$ws.Range("F33").Value = 25
$ws.Range("C34").Value = 33
$ws.Range("D34").Value = 'Nabumet 750mg FC Tab 24''s'
$ws.Range("E34").Value = '24''s'
$ws.Range("F34").Value = 132
$ws.Range("C35").Value = 34
$ws.Range("D35").Value = 'Nabumet 500mg FC Tab 30''s'
$ws.Range("E35").Value = '30''s'
$ws.Range("F35").Value = 1217
$ws.Range("C36").Value = 35
$ws.Range("D36").Value = 'Naprox 15gm Gel'
$ws.Range("E36").Value = '15 gm'
$ws.Range("F36").Value = 61
$ws.Range("C37").Value = 36
$ws.Range("D37").Value = 'Naprox 250mg Tablet'
$ws.Range("E37").Value = ' 50 ''s'
$ws.Range("F37").Value = 255
$ws.Range("C38").Value = 37
$ws.Range("D38").Value = 'Naprox 500mg Tablet'
$ws.Range("E38").Value = '50 ''s'
$ws.Range("F38").Value = 758
$ws.Range("C39").Value = 38
$ws.Range("D39").Value = 'Naprox 50ml Suspension'
$ws.Range("E39").Value = '50 ml'
$ws.Range("F39").Value = 386
$ws.Range("C40").Value = 39
$ws.Range("D40").Value = 'Naprox Plus 375mg Tablet - 30''s'
$ws.Range("E40").Value = '30 ''s'
$ws.Range("F40").Value = 373
$ws.Range("C41").Value = 40
$ws.Range("D41").Value = 'Naprox Plus 500mg Tablet - 36''s'
$ws.Range("E41").Value = '36''s'
$ws.Range("F41").Value = 802
$ws.Range("C42").Value = 41
$ws.Range("D42").Value = 'Ontin 60ml Syrup'
$ws.Range("E42").Value = '60 ml'
$ws.Range("F42").Value = 613
$ws.Range("C43").Value = 42
$ws.Range("D43").Value = 'Ontin 10mg Tablet'
$ws.Range("E43").Value = '100 ''s'
$ws.Range("F43").Value = 41
$ws.Range("C44").Value = 43
$ws.Range("D44").Value = 'Oradin 60ml Suspension'
$ws.Range("E44").Value = '60 ml'
$ws.Range("F44").Value = 1684
$ws.Range("C45").Value = 44
$ws.Range("D45").Value = 'Oradin FT 10mg Tablet'
$ws.Range("E45").Value = '40 ''s'
$ws.Range("F45").Value = 53
$ws.Range("C46").Value = 45
$ws.Range("D46").Value = 'Oradin 10mg Tablet'
$ws.Range("E46").Value = '100 ''s'
$ws.Range("F46").Value = 1316
$ws.Range("C47").Value = 46
$ws.Range("D47").Value = 'Osticare Tablet 30''s'
$ws.Range("E47").Value = '30''s'
$ws.Range("F47").Value = 330
$ws.Range("C48").Value = 47
$ws.Range("D48").Value = 'Osticare FC Tab Container 30''s'
$ws.Range("E48").Value = '30''s'
$ws.Range("F48").Value = 12
$ws.Range("C49").Value = 48
$ws.Range("D49").Value = 'Paino 100mg Tablet'
$ws.Range("E49").Value = '100 ''s'
$ws.Range("F49").Value = 168
$ws.Range("C50").Value = 49
$ws.Range("D50").Value = 'Quinox DS 60ml Pellets for Suspension'
$ws.Range("E50").Value = '60 ml'
$ws.Range("F50").Value = 1343
$ws.Range("C51").Value = 50
$ws.Range("D51").Value = 'Quinox 750mg Tablet'
$ws.Range("E51").Value = '10 ''s'
$ws.Range("F51").Value = 17
$ws.Range("C52").Value = 51
$ws.Range("D52").Value = 'Quinox 500mg Tablet - 20''s'
$ws.Range("E52").Value = '20 ''s'
$ws.Range("F52").Value = 50
$ws.Range("C53").Value = 52
$ws.Range("D53").Value = 'Quinox 500mg Tablet (40''s)'
$ws.Range("E53").Value = '40 ''s'
$ws.Range("F53").Value = 276
$ws.Range("C54").Value = 53
$ws.Range("D54").Value = 'Quinox 250mg Tablet'
$ws.Range("E54").Value = '30 ''s'
$ws.Range("F54").Value = 16
$ws.Range("C55").Value = 54
$ws.Range("D55").Value = 'Rupaday 10mg Tablet 30''s'
$ws.Range("E55").Value = '30''s'
$ws.Range("F55").Value = 151
$ws.Range("C56").Value = 55
$ws.Range("D56").Value = 'Sk-Mox 100ml Dry Suspension'
$ws.Range("E56").Value = '100 ml'
$ws.Range("F56").Value = 990
$ws.Range("C57").Value = 56
$ws.Range("D57").Value = 'Sk-Mox DS 100ml Dry Suspension'
$ws.Range("E57").Value = '100 ml'
$ws.Range("F57").Value = 82
$ws.Range("C58").Value = 57
$ws.Range("D58").Value = 'Sk-Mox 15ml P/D'
$ws.Range("E58").Value = '15 ml'
$ws.Range("F58").Value = 213
$ws.Range("C59").Value = 58
$ws.Range("D59").Value = 'Sk-Mox 500mg Capsule 50''s'
$ws.Range("E59").Value = '50''s'
$ws.Range("F59").Value = 1169
$ws.Range("C60").Value = 59
$ws.Range("D60").Value = 'Sk-Mox 250mg Capsule'
$ws.Range("E60").Value = '100 ''s'
$ws.Range("F60").Value = 53
$ws.Range("C61").Value = 60
$ws.Range("D61").Value = 'Sk-Mox 500mg Capsule'
$ws.Range("E61").Value = '48 ''s'
$ws.Range("F61").Value = 2
$ws.Range("C62").Value = 61
$ws.Range("D62").Value = 'Stiba 30ml Syrup'
$ws.Range("E62").Value = '30ml'
$ws.Range("F62").Value = 60
$ws.Range("C63").Value = 62
$ws.Range("D63").Value = 'Stiba 10mg Tablet - 30''s'
$ws.Range("E63").Value = '30''s'
$ws.Range("F63").Value = 183
$ws.Range("C64").Value = 63
$ws.Range("D64").Value = 'Sulidac 100mg Tablet 50''s'
$ws.Range("E64").Value = '50''s'
$ws.Range("F64").Value = 64
$ws.Range("C65").Value = 64
$ws.Range("D65").Value = 'Sulidac 200 Tablet 20''s'
$ws.Range("E65").Value = '20''s'
$ws.Range("F65").Value = 252
$ws.Range("C66").Value = 65
$ws.Range("D66").Value = 'Tenoxim 20mg Tablet'
$ws.Range("E66").Value = '30''s'
$ws.Range("F66").Value = 96
$ws.Range("C67").Value = 66
$ws.Range("D67").Value = 'Timothy 50mg Tablet'
$ws.Range("E67").Value = '50 ''s'
$ws.Range("F67").Value = 342
$ws.Range("C68").Value = 67
$ws.Range("D68").Value = 'Timothy 5mg IM/IV'
$ws.Range("E68").Value = '5''s'
$ws.Range("F68").Value = 98
$ws.Range("C69").Value = 68
$ws.Range("D69").Value = 'Tojak 5mg Tablet 10''s'
$ws.Range("E69").Value = '10''S'
$ws.Range("F69").Value = 278
$ws.Range("C70").Value = 69
$ws.Range("D70").Value = 'Toperin 50mg Tablet - 60''s'
$ws.Range("E70").Value = '60 ''s'
$ws.Range("F70").Value = 36
$ws.Range("C71").Value = 70
$ws.Range("D71").Value = 'Toti Tablet'
$ws.Range("E71").Value = '100 ''s'
$ws.Range("F71").Value = 194
$ws.Range("C72").Value = 71
$ws.Range("D72").Value = 'Toti 100ml Syrup'
$ws.Range("E72").Value = '100 ml'
$ws.Range("F72").Value = 8533
$ws.Range("C73").Value = 72
$ws.Range("D73").Value = 'Visomox 400mg FC Tablet 10''s'
$ws.Range("E73").Value = '10''s'
$ws.Range("F73").Value = 119
$ws.Range("C74").Value = 73
$ws.Range("D74").Value = 'Volmax SR 100mg Capsule'
$ws.Range("E74").Value = '60 ''s'
$ws.Range("F74").Value = 4
$ws.Range("C75").Value = 74
$ws.Range("D75").Value = 'Xenthol 30 Cream'
$ws.Range("E75").Value = '15gm'
$ws.Range("F75").Value = 64
$ws.Range("C76").Value = 75
$ws.Range("D76").Value = 'Zithrox 20ml Powder for Suspension'
$ws.Range("E76").Value = '20ml'
$ws.Range("F76").Value = 687
$ws.Range("C77").Value = 76
$ws.Range("D77").Value = 'Zithrox 35ml Dry Suspension'
$ws.Range("E77").Value = '35ml'
$ws.Range("F77").Value = 1678
$ws.Range("C78").Value = 77
$ws.Range("D78").Value = 'Zithrox 500mg Tablet - 12''s'
$ws.Range("E78").Value = '12 ''s'
$ws.Range("F78").Value = 1456
$ws.Range("C79").Value = 78
$ws.Range("D79").Value = 'Zithrox 250mg Tablet - 12''s'
$ws.Range("E79").Value = '12''s'
$ws.Range("F79").Value = 198
$ws.Range("C80").Value = 79
$ws.Range("D80").Value = 'Zithrox 50ml Powder for Suspension'
$ws.Range("E80").Value = '50ml'
$ws.Range("F80").Value = 192
